$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Activités")

# Row 23: fill in the new journal entry
$ws.Range("A23").Value = 44266
$ws.Range("B23").Value = 0.57430555555555551
$ws.Range("C23").Value = 0.62638888888888888
$ws.Range("E23").Value = "Rédaction documentation"
$ws.Range("G23").Value = "Avancement sur la documentation"

# Update the active selection to G23 (as captured in the sheet view)
$ws.Range("G23").Select()
